$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.601.83"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "3.525.69"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'610.75"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "'151.81"
$ws.Range("E6").Value = "  -1.43%  "
$ws.Range("D7").Value = "3.525.03"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("D11").Value = "'7.06"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "'0.426"
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("D13").Value = "'0.0000222"
$ws.Range("E13").Value = "  -1.55%  "
$ws.Range("D14").Value = "4.117.84"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "'32.01"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "3.524.80"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "67.506.67"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("D21").Value = "'446.78"
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("D22").Value = "'9.30"
$ws.Range("E22").Value = "  -4.45%  "
$ws.Range("D23").Value = "'0.625"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("D24").Value = "'77.40"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  +11.51%  "
$ws.Range("D26").Value = "3.667.08"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'10.22"
$ws.Range("E28").Value = "  -4.50%  "
$ws.Range("D29").Value = "'8.38"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("E31").Value = "  -3.78%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "'0.165"
$ws.Range("E33").Value = "  +4.66%  "
$ws.Range("D34").Value = "'25.80"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "'6.15"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("D36").Value = "3.517.21"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "  -3.39%  "
$ws.Range("D38").Value = "'8.07"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'177.28"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "'2.20"
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("D43").Value = "'0.0877"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("E44").Value = "  -3.34%  "
$ws.Range("D45").Value = "'0.882"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").Value = "'45.60"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").Value = "'2.63"
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("E48").Value = "  +5.42%  "
$ws.Range("D49").Value = "'27.18"
$ws.Range("E49").Value = "  -5.00%  "
$ws.Range("D50").Value = "'7.60"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -1.02%  "
